$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (sending cluster "ECs" group no longer present in new TPM data)
$ws.Rows("8:10").Delete()

# Update rows 2-7 with the new TPM-derived values
# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il27"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.148825
$ws.Range("H2").Value = 0.446475
$ws.Range("I2").Value = 0.2437842456871251
$ws.Range("J2").Value = 0.2437842456871251
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.185458333333334
$ws.Range("N2").Value = 6.556375
$ws.Range("O2").Value = 0.4976439168256567
$ws.Range("P2").Value = 0.4976439168256567
$ws.Range("Q2").Value = 0.3252508364583334
$ws.Range("R2").Value = 2.927257528125
$ws.Range("S2").Value = 0.1213177468841292
$ws.Range("T2").Value = 0.1213177468841292

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il27"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.148825
$ws.Range("H3").Value = 0.446475
$ws.Range("I3").Value = 0.2437842456871251
$ws.Range("J3").Value = 0.2437842456871251
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.742815333333333
$ws.Range("N3").Value = 5.228446
$ws.Range("O3").Value = 0.396851056620684
$ws.Range("P3").Value = 0.396851056620684
$ws.Range("Q3").Value = 0.2593744919833333
$ws.Range("R3").Value = 2.33437042785
$ws.Range("S3").Value = 0.09674603548841205
$ws.Range("T3").Value = 0.09674603548841205

# Row 4: FAPs -> Resolving-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Il27"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.148825
$ws.Range("H4").Value = 0.446475
$ws.Range("I4").Value = 0.2437842456871251
$ws.Range("J4").Value = 0.2437842456871251
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4633369999999999
$ws.Range("N4").Value = 1.390011
$ws.Range("O4").Value = 0.1055050265536593
$ws.Range("P4").Value = 0.1055050265536593
$ws.Range("Q4").Value = 0.068956129025
$ws.Range("R4").Value = 0.620605161225
$ws.Range("S4").Value = 0.02572046331458393
$ws.Range("T4").Value = 0.02572046331458393

# Row 5: Resolving-Mac -> ECs
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Il27"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4616533333333333
$ws.Range("H5").Value = 1.38496
$ws.Range("I5").Value = 0.7562157543128748
$ws.Range("J5").Value = 0.7562157543128749
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.185458333333334
$ws.Range("N5").Value = 6.556375
$ws.Range("O5").Value = 0.4976439168256567
$ws.Range("P5").Value = 0.4976439168256567
$ws.Range("Q5").Value = 1.008924124444444
$ws.Range("R5").Value = 9.08031712
$ws.Range("S5").Value = 0.3763261699415275
$ws.Range("T5").Value = 0.3763261699415276

# Row 6: Resolving-Mac -> FAPs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Il27"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4616533333333333
$ws.Range("H6").Value = 1.38496
$ws.Range("I6").Value = 0.7562157543128748
$ws.Range("J6").Value = 0.7562157543128749
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.742815333333333
$ws.Range("N6").Value = 5.228446
$ws.Range("O6").Value = 0.396851056620684
$ws.Range("P6").Value = 0.396851056620684
$ws.Range("Q6").Value = 0.8045765080177778
$ws.Range("R6").Value = 7.24118857216
$ws.Range("S6").Value = 0.3001050211322719
$ws.Range("T6").Value = 0.300105021132272

# Row 7: Resolving-Mac -> Resolving-Mac
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Il27"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4616533333333333
$ws.Range("H7").Value = 1.38496
$ws.Range("I7").Value = 0.7562157543128748
$ws.Range("J7").Value = 0.7562157543128749
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4633369999999999
$ws.Range("N7").Value = 1.390011
$ws.Range("O7").Value = 0.1055050265536593
$ws.Range("P7").Value = 0.1055050265536593
$ws.Range("Q7").Value = 0.2139010705066666
$ws.Range("R7").Value = 1.92510963456
$ws.Range("S7").Value = 0.07978456323907532
$ws.Range("T7").Value = 0.07978456323907535
